$wb = $excel.ActiveWorkbook

# Add the new "data2021" worksheet as the first sheet in the workbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "data2021"
$newSheet.Move($wb.Worksheets.Item(1))

$ws = $wb.Worksheets.Item("data2021")

$ws.Cells.Item(1,1).Value = "woj."
$ws.Cells.Item(1,2).Value = "area"
$ws.Cells.Item(1,3).Value = "co2"
$ws.Cells.Item(1,4).Value = "metan"
$ws.Cells.Item(1,5).Value = "n2o"
$ws.Cells.Item(1,6).Value = "so2"
$ws.Cells.Item(1,7).Value = "no"
$ws.Cells.Item(1,8).Value = "co"
$ws.Cells.Item(2,1).Value = "dolnoslaskie"
$ws.Cells.Item(2,2).Value = 19947
$ws.Cells.Item(2,3).Value = 24306.62
$ws.Cells.Item(2,4).Value = 42.2
$ws.Cells.Item(2,5).Value = 3.18
$ws.Cells.Item(2,6).Value = 26.56
$ws.Cells.Item(2,7).Value = 37.65
$ws.Cells.Item(2,8).Value = 158.11
$ws.Cells.Item(3,1).Value = "kujawsko-pomorskie"
$ws.Cells.Item(3,2).Value = 17972
$ws.Cells.Item(3,3).Value = 15750.88
$ws.Cells.Item(3,4).Value = 62.3
$ws.Cells.Item(3,5).Value = 6.2
$ws.Cells.Item(3,6).Value = 27.18
$ws.Cells.Item(3,7).Value = 44.77
$ws.Cells.Item(3,8).Value = 175.51
$ws.Cells.Item(4,1).Value = "lubelskie"
$ws.Cells.Item(4,2).Value = 25122
$ws.Cells.Item(4,3).Value = 12270.73
$ws.Cells.Item(4,4).Value = 130.4
$ws.Cells.Item(4,5).Value = 6.11
$ws.Cells.Item(4,6).Value = 19.91
$ws.Cells.Item(4,7).Value = 31.41
$ws.Cells.Item(4,8).Value = 145.78
$ws.Cells.Item(5,1).Value = "lubuskie"
$ws.Cells.Item(5,2).Value = 13988
$ws.Cells.Item(5,3).Value = 4986.79
$ws.Cells.Item(5,4).Value = 40.2
$ws.Cells.Item(5,5).Value = 1.68
$ws.Cells.Item(5,6).Value = 3.98
$ws.Cells.Item(5,7).Value = 12.54
$ws.Cells.Item(5,8).Value = 61.75
$ws.Cells.Item(6,1).Value = "lodzkie"
$ws.Cells.Item(6,2).Value = 18219
$ws.Cells.Item(6,3).Value = 48236.01
$ws.Cells.Item(6,4).Value = 93.93
$ws.Cells.Item(6,5).Value = 5.26
$ws.Cells.Item(6,6).Value = 67.25
$ws.Cells.Item(6,7).Value = 66.62
$ws.Cells.Item(6,8).Value = 229.63
$ws.Cells.Item(7,1).Value = "malopolskie"
$ws.Cells.Item(7,2).Value = 15183
$ws.Cells.Item(7,3).Value = 17924.34
$ws.Cells.Item(7,4).Value = 45.04
$ws.Cells.Item(7,5).Value = 2.35
$ws.Cells.Item(7,6).Value = 22.78
$ws.Cells.Item(7,7).Value = 32.44
$ws.Cells.Item(7,8).Value = 156.95
$ws.Cells.Item(8,1).Value = "mazowieckie"
$ws.Cells.Item(8,2).Value = 35558
$ws.Cells.Item(8,3).Value = 58045.52
$ws.Cells.Item(8,4).Value = 148.77
$ws.Cells.Item(8,5).Value = 9.9
$ws.Cells.Item(8,6).Value = 49.15
$ws.Cells.Item(8,7).Value = 87.88
$ws.Cells.Item(8,8).Value = 229.28
$ws.Cells.Item(9,1).Value = "opolskie"
$ws.Cells.Item(9,2).Value = 9412
$ws.Cells.Item(9,3).Value = 20797.11
$ws.Cells.Item(9,4).Value = 23.12
$ws.Cells.Item(9,5).Value = 2.76
$ws.Cells.Item(9,6).Value = 13.43
$ws.Cells.Item(9,7).Value = 25.71
$ws.Cells.Item(9,8).Value = 68.32
$ws.Cells.Item(10,1).Value = "podkarpackie"
$ws.Cells.Item(10,2).Value = 17846
$ws.Cells.Item(10,3).Value = 8028.31
$ws.Cells.Item(10,4).Value = 31.42
$ws.Cells.Item(10,5).Value = 1.54
$ws.Cells.Item(10,6).Value = 16.39
$ws.Cells.Item(10,7).Value = 19.26
$ws.Cells.Item(10,8).Value = 242.66
$ws.Cells.Item(11,1).Value = "podlaskie"
$ws.Cells.Item(11,2).Value = 20187
$ws.Cells.Item(11,3).Value = 4119.51
$ws.Cells.Item(11,4).Value = 104.07
$ws.Cells.Item(11,5).Value = 7.32
$ws.Cells.Item(11,6).Value = "string"
$ws.Cells.Item(11,7).Value = 17.6
$ws.Cells.Item(11,8).Value = 68.64
$ws.Cells.Item(12,1).Value = "pomorskie"
$ws.Cells.Item(12,2).Value = 18310
$ws.Cells.Item(12,3).Value = 13564.65
$ws.Cells.Item(12,4).Value = 53.3
$ws.Cells.Item(12,5).Value = 3.92
$ws.Cells.Item(12,6).Value = 23.21
$ws.Cells.Item(12,7).Value = 35
$ws.Cells.Item(12,8).Value = 161.24
$ws.Cells.Item(13,1).Value = "slaskie"
$ws.Cells.Item(13,2).Value = 12333
$ws.Cells.Item(13,3).Value = 50417.5
$ws.Cells.Item(13,4).Value = 488.88
$ws.Cells.Item(13,5).Value = 2.61
$ws.Cells.Item(13,6).Value = 44.31
$ws.Cells.Item(13,7).Value = 55.84
$ws.Cells.Item(13,8).Value = 219.66
$ws.Cells.Item(14,1).Value = "swietokrzyskie"
$ws.Cells.Item(14,2).Value = 11711
$ws.Cells.Item(14,3).Value = 15333.41
$ws.Cells.Item(14,4).Value = 22.97
$ws.Cells.Item(14,5).Value = 1.87
$ws.Cells.Item(14,6).Value = 16.02
$ws.Cells.Item(14,7).Value = 22.2
$ws.Cells.Item(14,8).Value = 68.98
$ws.Cells.Item(15,1).Value = "warminsko-mazurskie"
$ws.Cells.Item(15,2).Value = 24173
$ws.Cells.Item(15,3).Value = 5524.27
$ws.Cells.Item(15,4).Value = 58.35
$ws.Cells.Item(15,5).Value = 5.67
$ws.Cells.Item(15,6).Value = 12.09
$ws.Cells.Item(15,7).Value = 21.25
$ws.Cells.Item(15,8).Value = 124.8
$ws.Cells.Item(16,1).Value = "wielkopolskie"
$ws.Cells.Item(16,2).Value = 29826
$ws.Cells.Item(16,3).Value = 21070.68
$ws.Cells.Item(16,4).Value = 141.33
$ws.Cells.Item(16,5).Value = 11.73
$ws.Cells.Item(16,6).Value = 23.5
$ws.Cells.Item(16,7).Value = 50.53
$ws.Cells.Item(16,8).Value = 213.79
$ws.Cells.Item(17,1).Value = "zachodniopomorskie"
$ws.Cells.Item(17,2).Value = 22892
$ws.Cells.Item(17,3).Value = 11199.94
$ws.Cells.Item(17,4).Value = 38.89
$ws.Cells.Item(17,5).Value = 5.69
$ws.Cells.Item(17,6).Value = 19.47
$ws.Cells.Item(17,7).Value = 30.68
$ws.Cells.Item(17,8).Value = 195.65

# Select F11 on the new sheet to match the saved selection
$ws.Range("F11").Select()

# Restore the active/selected sheet and its previous selection on the old "data2020" sheet
$ws2020 = $wb.Worksheets.Item("data2020")
$ws2020.Range("H5").Select()

$ws.Select()

Write-Host "done"
